# This workbook is an AHB-Diff export with header row:
#   A..J  = "<Label>_old"   -> rename to "<Label>_FV2404"
#   K     = "diff"          -> unchanged
#   L..U  = "<Label>_new"   -> rename to "<Label>_FV2410"
# After renaming we turn the data range into a proper Excel Table
# ("Table1") and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$namesFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$namesFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# 1. Rename the header row cells: columns A-J (formerly "_old") and
#    columns L-U (formerly "_new"). Column K ("diff") is left as-is.
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $namesFV2404[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $namesFV2410[$i]
}

# 2. Turn the whole used range (header + data) into an Excel Table.
$range = $ws.Range("A1:U74")
$listObject = $ws.ListObjects.Add(1, $range, [System.Type]::Missing, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = $null

# 3. Freeze the header row (split below row 1, top-left of the
#    scrollable area is A2).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
